$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A55").Value = "2025-04-29 07:00:27"
$ws.Range("B55").Value = 146
